# "check in and undo check in"
# Rename the "edit_booking_1_player_id" column on the "Test" sheet to
# "check_in_bag_id" with values CI_001..CI_005, and add a new column
# "undo_check_in_bag_id" with value "UCI_001" for every data row.
# Also move the active tab selection from "booking_at_tee_time" to "Test".

$wb = $excel.ActiveWorkbook

$testSheet = $wb.Worksheets.Item("Test")

# Rename header H1/I1 and fill in the new check-in / undo-check-in ids,
# writing row by row (H then I) so new shared strings are interned in the
# same order the original author's edit produced them in.
$testSheet.Range("H1").Value = "check_in_bag_id"
$testSheet.Range("I1").Value = "undo_check_in_bag_id"

$testSheet.Range("H2").Value = "CI_001"
$testSheet.Range("I2").Value = "UCI_001"

$testSheet.Range("H3").Value = "CI_002"
$testSheet.Range("I3").Value = "UCI_001"

$testSheet.Range("H4").Value = "CI_003"
$testSheet.Range("I4").Value = "UCI_001"

$testSheet.Range("H5").Value = "CI_004"
$testSheet.Range("I5").Value = "UCI_001"

$testSheet.Range("H6").Value = "CI_005"
$testSheet.Range("I6").Value = "UCI_001"

# Copy the style used by column H onto the new column I so it matches formatting.
$testSheet.Range("H1").Copy() | Out-Null
$testSheet.Range("I1").PasteSpecial(-4122) | Out-Null
$testSheet.Range("H2:H6").Copy() | Out-Null
$testSheet.Range("I2:I6").PasteSpecial(-4122) | Out-Null

# Set a sensible column width for the new/changed columns (nearest width
# achievable through the ColumnWidth setter's pixel-rounded granularity).
$testSheet.Columns.Item(8).ColumnWidth = 18.5
$testSheet.Columns.Item(9).ColumnWidth = 24.16

# Update the selection on flow_4_player first (selecting a range on another
# sheet activates it as a side effect, so this must not be the last action).
$wb.Worksheets.Item("flow_4_player").Range("E21").Select() | Out-Null

# Finally, move the active tab from booking_at_tee_time to Test: activating
# Test and selecting its new cell must happen last so Test ends up as the
# workbook's active sheet (booking_at_tee_time keeps its own selection at
# C20, it just stops being the selected tab).
$testSheet.Activate()
$testSheet.Range("F22").Select() | Out-Null
